$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("!MonsterGroup")

# Rename header "Kind" -> "Index"
$ws.Range("A1").Value = "Index"

# Renumber the Index column for rows 4..8 (1-based running index)
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# Update selection to A8 as in the saved file
$ws.Range("A8").Select()
